# Update the worksheet date and all 25 "two digit / one digit" division
# problems in the table to the next day's values.
#
# NOTE on ordering: each Find/replace below targets exactly one occurrence
# (Replace = 1 => wdReplaceOne), so the calls must run in an order where no
# earlier replacement's *new* text collides with a *later* search's target
# text. The only colliding pair here is "36÷5=" -> "19÷6=" together with the
# original "19÷6=" -> "66÷4=": the original "19÷6=" cell must be replaced
# first, before "36÷5=" turns into a fresh "19÷6=". All other replacements
# are independent, so their relative order doesn't matter.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-12-15 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-16 Monday", 1) | Out-Null
$d.Content.Find.Execute("20÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "91÷7=", 1) | Out-Null
$d.Content.Find.Execute("85÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "52÷9=", 1) | Out-Null
$d.Content.Find.Execute("71÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "46÷5=", 1) | Out-Null
$d.Content.Find.Execute("19÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "66÷4=", 1) | Out-Null
$d.Content.Find.Execute("36÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "19÷6=", 1) | Out-Null
$d.Content.Find.Execute("60÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "64÷9=", 1) | Out-Null
$d.Content.Find.Execute("96÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷6=", 1) | Out-Null
$d.Content.Find.Execute("41÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "52÷7=", 1) | Out-Null
$d.Content.Find.Execute("67÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "73÷8=", 1) | Out-Null
$d.Content.Find.Execute("52÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷6=", 1) | Out-Null
$d.Content.Find.Execute("44÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "60÷4=", 1) | Out-Null
$d.Content.Find.Execute("93÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "49÷7=", 1) | Out-Null
$d.Content.Find.Execute("14÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷9=", 1) | Out-Null
$d.Content.Find.Execute("49÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "27÷5=", 1) | Out-Null
$d.Content.Find.Execute("77÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "47÷9=", 1) | Out-Null
$d.Content.Find.Execute("46÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷2=", 1) | Out-Null
$d.Content.Find.Execute("25÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "83÷6=", 1) | Out-Null
$d.Content.Find.Execute("75÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "73÷4=", 1) | Out-Null
$d.Content.Find.Execute("61÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "84÷2=", 1) | Out-Null
$d.Content.Find.Execute("18÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "95÷4=", 1) | Out-Null
$d.Content.Find.Execute("89÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "33÷7=", 1) | Out-Null
$d.Content.Find.Execute("85÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "13÷9=", 1) | Out-Null
$d.Content.Find.Execute("94÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷5=", 1) | Out-Null
$d.Content.Find.Execute("80÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "22÷6=", 1) | Out-Null
$d.Content.Find.Execute("62÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷6=", 1) | Out-Null
